$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.199790954589844
$ws.Range("B1").Value = 2.063236951828003
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.083468198776245
$ws.Range("E1").Value = 1.20766007900238
